$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 1 (headers) - a new "Alcance" header is inserted at D1 and the rest of
# the header row shifts one column to the right (E1..M1); N1 becomes blank.
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Iteración"
$ws.Range("B1").Value = "Historia(s)"
$ws.Range("C1").Value = "Rama (feature/*)"
$ws.Range("D1").Value = "Alcance (conceptual/lógico/físico/creacion_de_historias)"
$ws.Range("E1").Value = "Autor"
$ws.Range("F1").Value = "Revisor"
$ws.Range("G1").Value = "PR URL"
$ws.Range("H1").Value = "Commit/Tag"
$ws.Range("I1").Value = "Cambios clave (resumen)"
$ws.Range("J1").Value = "Fecha inicio"
$ws.Range("K1").Value = "Fecha cierre"
$ws.Range("L1").Value = "Estado (Planificado/En curso/Completado)"
$ws.Range("M1").Value = "Archivos afectados (ruta)"
$ws.Range("N1").Value = ""

# ---------------------------------------------------------------------------
# Row 2 - single consolidated iteration entry describing the creation of the
# first 3 user stories.
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "1.0"
$ws.Range("B2").Value = "H1 – Bicicletas"
$ws.Range("C2").Value = "feature/creacion_historias_usuario"
$ws.Range("D2").Value = "creacion_de_historias"
$ws.Range("E2").Value = "andr4f"
$ws.Range("F2").Value = "Yineth Avila, Angel Trillo, Andres Penagos"
$ws.Range("G2").Value = "https://github.com/andr4f/bici-go-bd/pull/3"
$ws.Range("H2").Value = "v0.1-iter1-creacion"
$ws.Range("I2").Value = "Creacion primeras 3 historias de usuario a partir de H1 principal"
$ws.Range("J2").Value = [DateTime]"2025-10-03"
$ws.Range("J2").NumberFormat = "mm-dd-yy"
$ws.Range("K2").Value = [DateTime]"2025-10-03"
# Copy J2's format (and only J2's) onto K2 so that both cells share the very
# same cell style entry instead of the engine minting two near identical
# (but distinct) style records.
$ws.Range("J2").Copy()
$ws.Range("K2").PasteSpecial(-4122)
$ws.Range("L2").Value = "Completado"
$ws.Range("M2").Value = "documentacion/back_logs_historias_usuario/sub_historias_h1/primeras_tres_sub_historias.pdf"
$ws.Range("N2").ClearContents()

# ---------------------------------------------------------------------------
# Row 3 - the second sample row from the template is removed entirely; only
# the (now empty) hyperlink-styled G3 cell remains.
# ---------------------------------------------------------------------------
$ws.Range("A3:N3").ClearContents()

# ---------------------------------------------------------------------------
# Hyperlinks - drop the two demo hyperlinks and add a single one for the new
# G2 PR link.
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("G2"), "https://github.com/andr4f/bici-go-bd/pull/3")

# ---------------------------------------------------------------------------
# Column widths
# ---------------------------------------------------------------------------
$ws.Columns.Item(4).ColumnWidth = 52
$ws.Columns.Item(13).ColumnWidth = 80.16666666666667
$ws.Columns.Item(14).ColumnWidth = 77.66666666666667

# ---------------------------------------------------------------------------
# Selection / view
# ---------------------------------------------------------------------------
$ws.Range("G3").Select()

Write-Host "edit complete"
